$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4477.625
$ws.Range("I62").Value = 3403.2
$ws.Range("J62").Value = 6268.3335
$ws.Range("K62").Value = 3403.2
$ws.Range("L62").Value = 6268.3335
$ws.Range("M62").Value = -2779.2
$ws.Range("N62").Value = -7516.3335
$ws.Range("H65").Value = 4477.625
$ws.Range("I65").Value = 3403.2
$ws.Range("J65").Value = 6268.3335
$ws.Range("K65").Value = 17016
$ws.Range("L65").Value = 31341.6675
$ws.Range("M65").Value = -13896
$ws.Range("N65").Value = -37581.6675
$ws.Range("H86").Value = 1777.1072
$ws.Range("I86").Value = 1911
$ws.Range("J86").Value = 1494.4445
$ws.Range("K86").Value = 1911
$ws.Range("L86").Value = 1494.4445
$ws.Range("M86").Value = -788
$ws.Range("N86").Value = -3740.4445
$ws.Range("H89").Value = 1777.1072
$ws.Range("I89").Value = 1911
$ws.Range("J89").Value = 1494.4445
$ws.Range("K89").Value = 9555
$ws.Range("L89").Value = 7472.2225
$ws.Range("M89").Value = -3939
$ws.Range("N89").Value = -18704.2225
$ws.Range("H137").Value = 2504260.8
$ws.Range("I137").Value = 3708317.5
$ws.Range("J137").Value = 3527.3845
$ws.Range("K137").Value = 11124952.5
$ws.Range("L137").Value = 10582.1535
$ws.Range("M137").Value = -11122402.5
$ws.Range("N137").Value = -15682.1535
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 90
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 22
$ws.Range("N5").Value = -324
$ws.Range("H44").Value = 21949.666
$ws.Range("J44").Value = 21949.666
$ws.Range("L44").Value = 21949.666
$ws.Range("N44").Value = -22925.666
$ws.Range("H63").Value = 2800.1667
$ws.Range("I63").Value = 1125.25
$ws.Range("K63").Value = 1125.25
$ws.Range("M63").Value = -439.25
$ws.Range("H66").Value = 2800.1667
$ws.Range("I66").Value = 1125.25
$ws.Range("K66").Value = 5626.25
$ws.Range("M66").Value = -2194.25
$ws.Range("H74").Value = 1609.2858
$ws.Range("I74").Value = 1118.1724
$ws.Range("J74").Value = 3983
$ws.Range("K74").Value = 1118.1724
$ws.Range("L74").Value = 3983
$ws.Range("M74").Value = -244.1723999999999
$ws.Range("N74").Value = -5731
$ws.Range("H77").Value = 1609.2858
$ws.Range("I77").Value = 1118.1724
$ws.Range("J77").Value = 3983
$ws.Range("K77").Value = 5590.861999999999
$ws.Range("L77").Value = 19915
$ws.Range("M77").Value = -1222.861999999999
$ws.Range("N77").Value = -28651
$ws.Range("H80").Value = 29105.666
$ws.Range("J80").Value = 29105.666
$ws.Range("L80").Value = 29105.666
$ws.Range("N80").Value = -31101.666
$ws.Range("H83").Value = 29105.666
$ws.Range("J83").Value = 29105.666
$ws.Range("L83").Value = 87316.99800000001
$ws.Range("N83").Value = -97300.99800000001
$ws.Range("H122").Value = 3447.5334
$ws.Range("I122").Value = 2316.5
$ws.Range("K122").Value = 6949.5
$ws.Range("M122").Value = -4499.5
$ws.Range("H132").Value = 2232.6
$ws.Range("I132").Value = 1584.8
$ws.Range("J132").Value = 3528.2
$ws.Range("K132").Value = 4754.4
$ws.Range("L132").Value = 10584.6
$ws.Range("M132").Value = -2224.4
$ws.Range("N132").Value = -15644.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 90
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = -330
$ws.Range("H38").Value = 62036
$ws.Range("J38").Value = 62036
$ws.Range("L38").Value = 62036
$ws.Range("N38").Value = -62868
$ws.Range("H82").Value = 21062.777
$ws.Range("J82").Value = 29127.666
$ws.Range("L82").Value = 29127.666
$ws.Range("N82").Value = -29893.666
$ws.Range("H85").Value = 21062.777
$ws.Range("J85").Value = 29127.666
$ws.Range("L85").Value = 29127.666
$ws.Range("N85").Value = -31779.666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 24092
$ws.Range("J50").Value = 24092
$ws.Range("L50").Value = 24092
$ws.Range("N50").Value = -25342
$ws.Range("H99").Value = 5210.875
$ws.Range("I99").Value = 3778.8333
$ws.Range("K99").Value = 3778.8333
$ws.Range("M99").Value = -2280.8333
$ws.Range("H122").Value = 2736
$ws.Range("I122").Value = 2407.6956
$ws.Range("J122").Value = 3814.7144
$ws.Range("K122").Value = 7223.0868
$ws.Range("L122").Value = 11444.1432
$ws.Range("M122").Value = -4773.0868
$ws.Range("N122").Value = -16344.1432
$ws.Range("H126").Value = 5210.875
$ws.Range("I126").Value = 3778.8333
$ws.Range("K126").Value = 11336.4999
$ws.Range("M126").Value = -8866.499899999999
$ws.Range("H134").Value = 2475.5476
$ws.Range("I134").Value = 1692.9062
$ws.Range("J134").Value = 4980
$ws.Range("K134").Value = 5078.7186
$ws.Range("L134").Value = 14940
$ws.Range("M134").Value = -2543.7186
$ws.Range("N134").Value = -20010
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 944.2727
$ws.Range("I107").Value = 485.33334
$ws.Range("J107").Value = 1262
$ws.Range("K107").Value = 1456.00002
$ws.Range("L107").Value = 3786
$ws.Range("M107").Value = 463.9999800000001
$ws.Range("N107").Value = -7626
$ws.Range("H117").Value = 2086.6924
$ws.Range("I117").Value = 799.75
$ws.Range("J117").Value = 2658.6667
$ws.Range("K117").Value = 2399.25
$ws.Range("L117").Value = 7976.000100000001
$ws.Range("M117").Value = 1042.75
$ws.Range("N117").Value = -14860.0001
$ws.Range("H131").Value = 1516.9803
$ws.Range("I131").Value = 3499.2856
$ws.Range("J131").Value = 1201.6136
$ws.Range("K131").Value = 10497.8568
$ws.Range("L131").Value = 3604.8408
$ws.Range("M131").Value = -5457.856800000001
$ws.Range("N131").Value = -13684.8408
$ws.Range("H136").Value = 2162.682
$ws.Range("I136").Value = 1638.6
$ws.Range("J136").Value = 3285.7144
$ws.Range("K136").Value = 4915.799999999999
$ws.Range("L136").Value = 9857.143199999999
$ws.Range("M136").Value = 184.2000000000007
$ws.Range("N136").Value = -20057.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1747.9286
$ws.Range("I97").Value = 1150
$ws.Range("J97").Value = 2824.2
$ws.Range("K97").Value = 1150
$ws.Range("L97").Value = 2824.2
$ws.Range("M97").Value = -654
$ws.Range("N97").Value = -3816.2
$ws.Range("H122").Value = 7662.643
$ws.Range("I122").Value = 1384.5
$ws.Range("K122").Value = 4153.5
$ws.Range("M122").Value = -1703.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12496.667
$ws.Range("I40").Value = 16660
$ws.Range("J40").Value = 8333.333000000001
$ws.Range("K40").Value = 16660
$ws.Range("L40").Value = 8333.333000000001
$ws.Range("M40").Value = -16524
$ws.Range("N40").Value = -8605.333000000001
$ws.Range("H51").Value = 60084
$ws.Range("J51").Value = 60084
$ws.Range("L51").Value = 60084
$ws.Range("N51").Value = -61040
$ws.Range("H100").Value = 3095
$ws.Range("I100").Value = 1681.6666
$ws.Range("J100").Value = 3801.6667
$ws.Range("K100").Value = 1681.6666
$ws.Range("L100").Value = 3801.6667
$ws.Range("M100").Value = -1140.6666
$ws.Range("N100").Value = -4883.6667
$ws.Range("H122").Value = 5325
$ws.Range("I122").Value = 4220
$ws.Range("J122").Value = 7166.6665
$ws.Range("K122").Value = 12660
$ws.Range("L122").Value = 21499.9995
$ws.Range("M122").Value = -10210
$ws.Range("N122").Value = -26399.9995
$ws.Range("H132").Value = 3063.842
$ws.Range("I132").Value = 2321.3
$ws.Range("J132").Value = 3888.889
$ws.Range("K132").Value = 6963.900000000001
$ws.Range("L132").Value = 11666.667
$ws.Range("M132").Value = -4433.900000000001
$ws.Range("N132").Value = -16726.667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1468.72
$ws.Range("I122").Value = 1034.238
$ws.Range("K122").Value = 3102.714
$ws.Range("M122").Value = -652.7139999999999
$ws.Range("H136").Value = 2132.7273
$ws.Range("I136").Value = 1629.725
$ws.Range("J136").Value = 3474.0667
$ws.Range("K136").Value = 4889.174999999999
$ws.Range("L136").Value = 10422.2001
$ws.Range("M136").Value = -2339.174999999999
$ws.Range("N136").Value = -15522.2001
